# Repull data / push all data / mean calculation
# Update column F ("dSF") values for the rows whose data changed on re-pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 2
    6  = -4
    12 = -3
    21 = -2
    28 = -1
    29 = 1
    30 = -1
    31 = 0
    36 = 2
    38 = 1
    41 = 3
    42 = 0
    43 = 1
    54 = 1
    55 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
